$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 306, shifting existing rows 306:380 down to 307:381
$ws.Rows("306:306").Insert()

# Populate the newly inserted row with the new price-record data
$ws.Range("A306").Value = 5
$ws.Range("B306").Value = "Macroferia Regional de Talca"
$ws.Range("C306").Value = "Maule"
$ws.Range("D306").Value = 44754
$ws.Range("E306").Value = 7
$ws.Range("F306").Value = 100112032
$ws.Range("G306").Value = "Zapallo italiano"
$ws.Range("H306").Value = "Sin especificar"
$ws.Range("I306").Value = "Primera"
$ws.Range("J306").Value = 400
$ws.Range("K306").Value = 12000
$ws.Range("L306").Value = 12000
$ws.Range("M306").Value = 12000
$ws.Range("N306").Value = "`$/caja 50 unidades"
$ws.Range("O306").Value = "Región de Arica y Parinacota"
$ws.Range("P306").Value = 240
$ws.Range("Q306").Value = 50
$ws.Range("R306").Value = "Hortaliza"
